$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 40.34291466666667
$ws.Range("H2").Value = 121.028744
$ws.Range("I2").Value = 0.3404392602027053
$ws.Range("J2").Value = 0.3404392602027053
$ws.Range("M2").Value = 168.1098273333333
$ws.Range("N2").Value = 504.329482
$ws.Range("O2").Value = 0.2984182258032519
$ws.Range("P2").Value = 0.298418225803252
$ws.Range("Q2").Value = 6782.040418736734
$ws.Range("R2").Value = 61038.36376863061
$ws.Range("S2").Value = 0.1015932800234629
$ws.Range("T2").Value = 0.101593280023463
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 40.34291466666667
$ws.Range("H3").Value = 121.028744
$ws.Range("I3").Value = 0.3404392602027053
$ws.Range("J3").Value = 0.3404392602027053
$ws.Range("O3").Value = 0.2893586437755394
$ws.Range("P3").Value = 0.2893586437755394
$ws.Range("Q3").Value = 6576.14665563489
$ws.Range("R3").Value = 59185.31990071401
$ws.Range("S3").Value = 0.09850904262020277
$ws.Range("T3").Value = 0.09850904262020278
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 40.34291466666667
$ws.Range("H4").Value = 121.028744
$ws.Range("I4").Value = 0.3404392602027053
$ws.Range("J4").Value = 0.3404392602027053
$ws.Range("M4").Value = 165.99353
$ws.Range("N4").Value = 497.98059
$ws.Range("O4").Value = 0.294661504941043
$ws.Range("P4").Value = 0.294661504941043
$ws.Range("Q4").Value = 6696.662816008773
$ws.Range("R4").Value = 60269.96534407896
$ws.Range("S4").Value = 0.1003143447523445
$ws.Range("T4").Value = 0.1003143447523445
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 40.34291466666667
$ws.Range("H5").Value = 121.028744
$ws.Range("I5").Value = 0.3404392602027053
$ws.Range("J5").Value = 0.3404392602027053
$ws.Range("M5").Value = 66.22673433333334
$ws.Range("N5").Value = 198.680203
$ws.Range("O5").Value = 0.1175616254801657
$ws.Range("P5").Value = 0.1175616254801657
$ws.Range("Q5").Value = 2671.779491861671
$ws.Range("R5").Value = 24046.01542675504
$ws.Range("S5").Value = 0.0400225928066951
$ws.Range("T5").Value = 0.04002259280669511
$ws.Range("I6").Value = 0.3497297648481489
$ws.Range("J6").Value = 0.3497297648481489
$ws.Range("M6").Value = 168.1098273333333
$ws.Range("N6").Value = 504.329482
$ws.Range("O6").Value = 0.2984182258032519
$ws.Range("P6").Value = 0.298418225803252
$ws.Range("Q6").Value = 6967.120652956321
$ws.Range("R6").Value = 62704.08587660689
$ws.Range("S6").Value = 0.1043657359365731
$ws.Range("T6").Value = 0.1043657359365731
$ws.Range("I7").Value = 0.3497297648481489
$ws.Range("J7").Value = 0.3497297648481489
$ws.Range("O7").Value = 0.2893586437755394
$ws.Range("P7").Value = 0.2893586437755394
$ws.Range("S7").Value = 0.1011973304443987
$ws.Range("T7").Value = 0.1011973304443987
$ws.Range("I8").Value = 0.3497297648481489
$ws.Range("J8").Value = 0.3497297648481489
$ws.Range("M8").Value = 165.99353
$ws.Range("N8").Value = 497.98059
$ws.Range("O8").Value = 0.294661504941043
$ws.Range("P8").Value = 0.294661504941043
$ws.Range("Q8").Value = 6879.41311620639
$ws.Range("R8").Value = 61914.71804585751
$ws.Range("S8").Value = 0.1030518988328326
$ws.Range("T8").Value = 0.1030518988328326
$ws.Range("I9").Value = 0.3497297648481489
$ws.Range("J9").Value = 0.3497297648481489
$ws.Range("M9").Value = 66.22673433333334
$ws.Range("N9").Value = 198.680203
$ws.Range("O9").Value = 0.1175616254801657
$ws.Range("P9").Value = 0.1175616254801657
$ws.Range("Q9").Value = 2744.691704648063
$ws.Range("R9").Value = 24702.22534183257
$ws.Range("S9").Value = 0.04111479963434449
$ws.Range("T9").Value = 0.04111479963434449
$ws.Range("G10").Value = 36.642055
$ws.Range("H10").Value = 109.926165
$ws.Range("I10").Value = 0.3092090445020276
$ws.Range("J10").Value = 0.3092090445020277
$ws.Range("M10").Value = 168.1098273333333
$ws.Range("N10").Value = 504.329482
$ws.Range("O10").Value = 0.2984182258032519
$ws.Range("P10").Value = 0.298418225803252
$ws.Range("Q10").Value = 6159.889539188503
$ws.Range("R10").Value = 55439.00585269653
$ws.Range("S10").Value = 0.09227361446261387
$ws.Range("T10").Value = 0.0922736144626139
$ws.Range("G11").Value = 36.642055
$ws.Range("H11").Value = 109.926165
$ws.Range("I11").Value = 0.3092090445020276
$ws.Range("J11").Value = 0.3092090445020277
$ws.Range("O11").Value = 0.2893586437755394
$ws.Range("P11").Value = 0.2893586437755394
$ws.Range("Q11").Value = 5972.883452640962
$ws.Range("R11").Value = 53755.95107376866
$ws.Range("S11").Value = 0.08947230976023714
$ws.Range("T11").Value = 0.08947230976023715
$ws.Range("G12").Value = 36.642055
$ws.Range("H12").Value = 109.926165
$ws.Range("I12").Value = 0.3092090445020276
$ws.Range("J12").Value = 0.3092090445020277
$ws.Range("M12").Value = 165.99353
$ws.Range("N12").Value = 497.98059
$ws.Range("O12").Value = 0.294661504941043
$ws.Range("P12").Value = 0.294661504941043
$ws.Range("Q12").Value = 6082.344055904149
$ws.Range("R12").Value = 54741.09650313735
$ws.Range("S12").Value = 0.0911120023943494
$ws.Range("T12").Value = 0.09111200239434941
$ws.Range("G13").Value = 36.642055
$ws.Range("H13").Value = 109.926165
$ws.Range("I13").Value = 0.3092090445020276
$ws.Range("J13").Value = 0.3092090445020277
$ws.Range("M13").Value = 66.22673433333334
$ws.Range("N13").Value = 198.680203
$ws.Range("O13").Value = 0.1175616254801657
$ws.Range("P13").Value = 0.1175616254801657
$ws.Range("Q13").Value = 2426.683641912388
$ws.Range("R13").Value = 21840.1527772115
$ws.Range("S13").Value = 0.03635111788482725
$ws.Range("T13").Value = 0.03635111788482726
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.07370033333333333
$ws.Range("H14").Value = 0.221101
$ws.Range("I14").Value = 0.0006219304471182344
$ws.Range("J14").Value = 0.0006219304471182345
$ws.Range("M14").Value = 168.1098273333333
$ws.Range("N14").Value = 504.329482
$ws.Range("O14").Value = 0.2984182258032519
$ws.Range("P14").Value = 0.298418225803252
$ws.Range("Q14").Value = 12.38975031107578
$ws.Range("R14").Value = 111.507752799682
$ws.Range("S14").Value = 0.0001855953806020467
$ws.Range("T14").Value = 0.0001855953806020468
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.07370033333333333
$ws.Range("H15").Value = 0.221101
$ws.Range("I15").Value = 0.0006219304471182344
$ws.Range("J15").Value = 0.0006219304471182345
$ws.Range("O15").Value = 0.2893586437755394
$ws.Range("P15").Value = 0.2893586437755394
$ws.Range("Q15").Value = 12.01361390404522
$ws.Range("R15").Value = 108.122525136407
$ws.Range("S15").Value = 0.0001799609507008472
$ws.Range("T15").Value = 0.0001799609507008472
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.07370033333333333
$ws.Range("H16").Value = 0.221101
$ws.Range("I16").Value = 0.0006219304471182344
$ws.Range("J16").Value = 0.0006219304471182345
$ws.Range("M16").Value = 165.99353
$ws.Range("N16").Value = 497.98059
$ws.Range("O16").Value = 0.294661504941043
$ws.Range("P16").Value = 0.294661504941043
$ws.Range("Q16").Value = 12.23377849217666
$ws.Range("R16").Value = 110.10400642959
$ws.Range("S16").Value = 0.0001832589615165147
$ws.Range("T16").Value = 0.0001832589615165147
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.07370033333333333
$ws.Range("H17").Value = 0.221101
$ws.Range("I17").Value = 0.0006219304471182344
$ws.Range("J17").Value = 0.0006219304471182345
$ws.Range("M17").Value = 66.22673433333334
$ws.Range("N17").Value = 198.680203
$ws.Range("O17").Value = 0.1175616254801657
$ws.Range("P17").Value = 0.1175616254801657
$ws.Range("Q17").Value = 2671.779491861671
$ws.Range("R17").Value = 24046.01542675504
$ws.Range("S17").Value = 0.0400225928066951
$ws.Range("T17").Value = 0.04002259280669511
